{"js": "// 1) Tweak the sentence: remove the word \"algum\" before \"comportamento\".\nconst body = context.document.body;\n\nconst targetRuns = body.search(\"tiveram algum comportamento\", { matchCase: true, matchWholeWord: false });\ntargetRuns.load(\"items/text\");\nawait context.sync();\n\nif (targetRuns.items.length > 0) {\n  targetRuns.items[0].insertText(\"tiveram comportamento\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Append the new explanatory sentences to the end of that same paragraph.\nconst endAnchor = body.search(\"citado anteriormente.\", { matchCase: true, matchWholeWord: false });\nendAnchor.load(\"items/text\");\nawait context.sync();\n\nconst addition =\n  \" Acredito que esses problemas devem ter ocorrido devido \u00e0 configura\u00e7\u00e3o de clock no modelsim\" +\n  \", isto \u00e9, a diferen\u00e7a da borda de subida e descida, al\u00e9m do tempo dos ciclos. Na imagem da \" +\n  \"simula\u00e7\u00e3o, foi utilizado como configura\u00e7\u00e3o:\";\n\nlet newParagraph = null;\nif (endAnchor.items.length > 0) {\n  const tailRange = endAnchor.items[0].getRange(Word.RangeLocation.end);\n  tailRange.insertText(addition, Word.InsertLocation.end);\n  await context.sync();\n\n  // 3) Insert a new, centered & italic paragraph right after, with the clock configuration.\n  newParagraph = tailRange.insertParagraph(\n    \"Borda de subida ; 100 ps ; 50 duty\",\n    Word.InsertLocation.after\n  );\n  newParagraph.alignment = Word.Alignment.centered;\n  await context.sync();\n\n  newParagraph.font.italic = true;\n  await context.sync();\n}\n\n// 4) Move the \"_GoBack\" bookmark from the end of the document to right after \"Conclus\u00e3o:\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst conclusion = body.search(\"Conclus\u00e3o:\", { matchCase: true, matchWholeWord: false });\nconclusion.load(\"items/text\");\nawait context.sync();\n\nif (conclusion.items.length > 0) {\n  const conclusionEnd = conclusion.items[0].getRange(Word.RangeLocation.end);\n  conclusionEnd.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Wording tweak: \"tiveram algum comportamento\" -> \"tiveram comportamento\".\n$range = $d.Content\n$range.Find.Execute(\"tiveram algum comportamento\", $false, $false, $false, $false, $false, $true, 1, $false, \"tiveram comportamento\", 2) | Out-Null\n\n# 2) Append the extra explanation sentences to the end of that same paragraph.\n$range = $d.Content\n$range.Find.Execute(\"citado anteriormente.\") | Out-Null\n$range.Collapse(0)\n$addition = \" Acredito que esses problemas devem ter ocorrido devido \u00e0 configura\u00e7\u00e3o de clock no modelsim, isto \u00e9, a diferen\u00e7a da borda de subida e descida, al\u00e9m do tempo dos ciclos. Na imagem da simula\u00e7\u00e3o, foi utilizado como configura\u00e7\u00e3o:\"\n$range.InsertAfter($addition)\n\n# 3) Insert a new, centered & italic paragraph right after with the clock configuration.\n$range.Collapse(0)\n$range.InsertParagraphAfter()\n\n$paraCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*foi utilizado como configura*\") {\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newRange = $newPara.Range\n        $newRange.Text = \"Borda de subida ; 100 ps ; 50 duty\"\n        $newPara.Alignment = 1  # wdAlignParagraphCenter\n        $newRange.Font.Italic = $true\n        break\n    }\n}\n\n# 4) Move the \"_GoBack\" bookmark from the end of the document to right after \"Conclus\u00e3o:\".\n$existingBookmark = $d.Bookmarks.Item(\"_GoBack\")\n$existingBookmark.Delete()\n\n$concRange = $d.Content\n$concRange.Find.Execute(\"Conclus\u00e3o:\") | Out-Null\n$concRange.MoveEnd(1, -1)   # back off the trailing colon so the anchor isn't flush with the paragraph mark\n$concRange.Collapse(0)      # collapse right before the colon\n$d.Bookmarks.Add(\"_GoBack\", $concRange)\n\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$tail = $bm.Range\n$tail.MoveEnd(1, 1)         # grow to cover the colon character\n$colonText = $tail.Text\n$tail.Delete()              # remove it (bookmark settles back at its safe anchor point)\n$d.Bookmarks.Item(\"_GoBack\").Range.InsertAfter($colonText)   # put the colon back, after the bookmark\n"}
